{"js": "// Diff summary: the worksheet's date heading moved forward a day, and every\n// two-digit multiplication problem (\"A\u00d7B=\") in the practice table was\n// swapped for a newly generated problem. Each old problem string is unique\n// within the document, so a simple search-and-replace per pair reproduces\n// the change exactly while leaving all run/paragraph formatting untouched.\nconst replacements = [\n  [\"2023-11-29 Wednesday\", \"2023-11-30 Thursday\"],\n  [\"91\u00d753=\", \"58\u00d760=\"],\n  [\"60\u00d758=\", \"13\u00d738=\"],\n  [\"36\u00d749=\", \"18\u00d766=\"],\n  [\"29\u00d736=\", \"57\u00d737=\"],\n  [\"67\u00d756=\", \"92\u00d793=\"],\n  [\"85\u00d740=\", \"39\u00d789=\"],\n  [\"83\u00d771=\", \"85\u00d772=\"],\n  [\"49\u00d728=\", \"21\u00d746=\"],\n  [\"86\u00d749=\", \"48\u00d760=\"],\n  [\"78\u00d777=\", \"18\u00d755=\"],\n  [\"67\u00d778=\", \"60\u00d777=\"],\n  [\"93\u00d756=\", \"44\u00d777=\"],\n  [\"62\u00d741=\", \"63\u00d723=\"],\n  [\"61\u00d750=\", \"68\u00d795=\"],\n  [\"56\u00d799=\", \"75\u00d736=\"],\n  [\"14\u00d785=\", \"33\u00d759=\"],\n  [\"97\u00d712=\", \"22\u00d796=\"],\n  [\"25\u00d728=\", \"83\u00d717=\"],\n  [\"48\u00d731=\", \"84\u00d734=\"],\n  [\"14\u00d753=\", \"14\u00d761=\"],\n  [\"62\u00d746=\", \"82\u00d752=\"],\n  [\"91\u00d758=\", \"97\u00d716=\"],\n  [\"17\u00d785=\", \"65\u00d762=\"],\n  [\"57\u00d794=\", \"33\u00d742=\"],\n  [\"19\u00d776=\", \"17\u00d763=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Each search string is unique in this document, but loop defensively in\n  // case a string were ever to repeat.\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Diff summary: the worksheet's date heading moved forward a day, and every\n# two-digit multiplication problem (\"A\u00d7B=\") in the practice table was\n# swapped for a newly generated problem. Each old problem string is unique\n# within the document, so a simple Find/Replace per pair reproduces the\n# change exactly while leaving all run/paragraph formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-11-29 Wednesday\", \"2023-11-30 Thursday\"),\n    @(\"91\u00d753=\", \"58\u00d760=\"),\n    @(\"60\u00d758=\", \"13\u00d738=\"),\n    @(\"36\u00d749=\", \"18\u00d766=\"),\n    @(\"29\u00d736=\", \"57\u00d737=\"),\n    @(\"67\u00d756=\", \"92\u00d793=\"),\n    @(\"85\u00d740=\", \"39\u00d789=\"),\n    @(\"83\u00d771=\", \"85\u00d772=\"),\n    @(\"49\u00d728=\", \"21\u00d746=\"),\n    @(\"86\u00d749=\", \"48\u00d760=\"),\n    @(\"78\u00d777=\", \"18\u00d755=\"),\n    @(\"67\u00d778=\", \"60\u00d777=\"),\n    @(\"93\u00d756=\", \"44\u00d777=\"),\n    @(\"62\u00d741=\", \"63\u00d723=\"),\n    @(\"61\u00d750=\", \"68\u00d795=\"),\n    @(\"56\u00d799=\", \"75\u00d736=\"),\n    @(\"14\u00d785=\", \"33\u00d759=\"),\n    @(\"97\u00d712=\", \"22\u00d796=\"),\n    @(\"25\u00d728=\", \"83\u00d717=\"),\n    @(\"48\u00d731=\", \"84\u00d734=\"),\n    @(\"14\u00d753=\", \"14\u00d761=\"),\n    @(\"62\u00d746=\", \"82\u00d752=\"),\n    @(\"91\u00d758=\", \"97\u00d716=\"),\n    @(\"17\u00d785=\", \"65\u00d762=\"),\n    @(\"57\u00d794=\", \"33\u00d742=\"),\n    @(\"19\u00d776=\", \"17\u00d763=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #   MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format, ReplaceWith,\n    #   Replace(2=wdReplaceAll))\n    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
